$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 166, shifting the existing rows 166-268 down to 167-269
$ws.Rows("166:166").Insert()

# Populate the newly inserted row 166 with its data
$ws.Range("A166").Value = 10
$ws.Range("B166").Value = "Vega Modelo de Temuco"
$ws.Range("C166").Value = "La Araucanía"
$ws.Range("D166").Value = 44777
$ws.Range("E166").Value = 9
$ws.Range("F166").Value = 100112039
$ws.Range("G166").Value = "Ciboulette"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 95
$ws.Range("K166").Value = 7000
$ws.Range("L166").Value = 7000
$ws.Range("M166").Value = 7000
$ws.Range("N166").Value = "$/docena de atados"
$ws.Range("O166").Value = "Provincia de Cautín"
$ws.Range("P166").Value = 2333
$ws.Range("Q166").Value = 3
$ws.Range("R166").Value = "Hortaliza"
